$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.510.17"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.075.08"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "234.53"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.01%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "608.57"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -2.07%  "
$ws.Range("E8").Value = "  -4.92%  "
$ws.Range("E9").Value = "  +0.08%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.801"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +6.10%  "
$ws.Range("D11").Value = "3.070.43"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "94.084.84"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  -4.53%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "33.61"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "3.642.72"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "3.041.57"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("E19").Value = "  -6.33%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "14.36"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.93%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.65"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "438.46"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.92%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "8.79"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.38%  "
$ws.Range("E24").Value = "  -7.25%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "8.34"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +5.88%  "
$ws.Range("E26").Value = "  -5.16%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "84.61"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.64%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "11.86"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "3.226.71"
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("E30").Value = "  +0.10%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.244"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +5.07%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.178"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +4.82%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.124"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -9.15%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "9.02"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.62%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "7.62"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -5.50%  "
$ws.Range("E36").Value = "  -3.08%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.891"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.37%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "25.32"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E39").Value = "  -1.81%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "24.01"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +3.85%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.435"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "468.23"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -4.68%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "3.69"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("E44").Value = "  -2.65%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.08"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -8.82%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "161.53"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.18%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.668"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.81"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.86%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "43.58"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("E51").Value = "  +0.10%  "
